# Refreshed examples and added CE quadrant counter
# - Adds three new transition parameter columns (parameter_5, parameter_6,
#   parameter_7) to the "transitions" sheet, ahead of the existing "notes"
#   column, mirroring the model-spec workbook refresh described in the
#   commit message.
# - Leaves a fresh-looking selection on each worksheet, then restores the
#   "utilities" sheet as the active tab (it was the active tab before the
#   edit and isn't touched by the content change).

$wb = $excel.ActiveWorkbook

# --- transitions sheet -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("transitions")

# Insert three blank columns right before the existing "notes" column (H)
# so it slides from H -> K, matching the target layout (A..G unchanged,
# H/I/J new, K = old H).
$null = $ws1.Columns("H:J").Insert()

$ws1.Range("H1").Value = "parameter_5"
$ws1.Range("I1").Value = "parameter_6"
$ws1.Range("J1").Value = "parameter_7"

$null = $ws1.Range("J1").Select()

# --- costs sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("costs")
$null = $ws2.Range("A2").Select()

# --- specification sheet ------------------------------------------------
$ws4 = $wb.Worksheets.Item("specification")
$null = $ws4.Range("A5").Select()

# --- restore the originally active sheet ("utilities") -------------------
$ws3 = $wb.Worksheets.Item("utilities")
$null = $ws3.Activate()
